# Update cryptocurrency price (D) and volume-change (E) columns with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new "Price" values are plain decimal numbers (e.g. "195.18"); Excel's COM layer
# would auto-detect and store these as numeric cells, but the source data keeps every
# Price/Volume cell as text. Force text storage for those specific cells via a temporary
# "@" (text) number format, then restore the default "Normal" style so no visible
# formatting change is introduced.
$forceTextRows = @(5, 6, 8, 9, 16, 19, 20, 21, 23, 24, 25, 28, 30, 32, 33, 35, 36, 37, 38, 40, 42, 43, 45, 46, 47, 49, 50, 51)
foreach ($r in $forceTextRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "76.300.86"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "2.859.23"
$ws.Range("E3").Value = "  +7.73%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "195.18"
$ws.Range("E5").Value = "  +5.00%  "
$ws.Range("D6").Value = "599.87"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  +3.97%  "
$ws.Range("D9").Value = "0.195"
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("D10").Value = "2.859.07"
$ws.Range("E10").Value = "  +7.87%  "
$ws.Range("E11").Value = "  +10.39%  "
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").Value = "3.389.06"
$ws.Range("E14").Value = "  +7.68%  "
$ws.Range("D15").Value = "76.069.76"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "27.54"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "2.856.33"
$ws.Range("E18").Value = "  +7.51%  "
$ws.Range("D19").Value = "9.08"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").Value = "12.47"
$ws.Range("E20").Value = "  +5.37%  "
$ws.Range("D21").Value = "383.09"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("E22").Value = "  +4.44%  "
$ws.Range("D23").Value = "4.15"
$ws.Range("E23").Value = "  +2.28%  "
$ws.Range("D24").Value = "72.00"
$ws.Range("E24").Value = "  +4.08%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "3.008.07"
$ws.Range("E26").Value = "  +7.71%  "
$ws.Range("E27").Value = "  +3.01%  "
$ws.Range("D28").Value = "9.75"
$ws.Range("E28").Value = "  +5.09%  "
$ws.Range("E29").Value = "  +12.58%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "514.62"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "7.73"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("E34").Value = "  +5.10%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "166.56"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Value = "20.04"
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "186.86"
$ws.Range("E40").Value = "  +10.52%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "0.346"
$ws.Range("E42").Value = "  +6.00%  "
$ws.Range("D43").Value = "5.10"
$ws.Range("E43").Value = "  +3.67%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").Value = "1.23"
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("D46").Value = "0.0899"
$ws.Range("E46").Value = "  +6.88%  "
$ws.Range("D47").Value = "40.32"
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "0.577"
$ws.Range("E49").Value = "  +10.25%  "
$ws.Range("D50").Value = "0.667"
$ws.Range("E50").Value = "  +13.46%  "
$ws.Range("D51").Value = "3.76"
$ws.Range("E51").Value = "  +4.10%  "

# Restore default styling on the cells we temporarily reformatted.
foreach ($r in $forceTextRows) {
    $ws.Range("D$r").Style = "Normal"
}

